$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, 8).Value = 76.28570999999999
$ws.Cells.Item(33, 9).Value = 76.28570999999999
$ws.Cells.Item(33, 11).Value = 76.28570999999999
$ws.Cells.Item(33, 13).Value = 152.71429
# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 2046.8928
$ws.Cells.Item(40, 9).Value = 1774.4375
$ws.Cells.Item(40, 10).Value = 2410.1667
$ws.Cells.Item(40, 11).Value = 1774.4375
$ws.Cells.Item(40, 12).Value = 2410.1667
$ws.Cells.Item(40, 13).Value = -1599.4375
$ws.Cells.Item(40, 14).Value = -2760.1667
# Row 64 (Leve Item ID 5506)
$ws.Cells.Item(64, 8).Value = 125004056
$ws.Cells.Item(64, 9).Value = 500001150
$ws.Cells.Item(64, 10).Value = 5026.6665
$ws.Cells.Item(64, 11).Value = 500001150
$ws.Cells.Item(64, 12).Value = 5026.6665
$ws.Cells.Item(64, 13).Value = -500000902
$ws.Cells.Item(64, 14).Value = -5522.6665
# Row 67 (Leve Item ID 5506)
$ws.Cells.Item(67, 8).Value = 125004056
$ws.Cells.Item(67, 9).Value = 500001150
$ws.Cells.Item(67, 10).Value = 5026.6665
$ws.Cells.Item(67, 11).Value = 500001150
$ws.Cells.Item(67, 12).Value = 5026.6665
$ws.Cells.Item(67, 13).Value = -500000292
$ws.Cells.Item(67, 14).Value = -6742.6665
# Row 69 (Leve Item ID 12616)
$ws.Cells.Item(69, 8).Value = 3902.6667
$ws.Cells.Item(69, 9).Value = 3925
$ws.Cells.Item(69, 10).Value = 3899.875
$ws.Cells.Item(69, 11).Value = 11775
$ws.Cells.Item(69, 12).Value = 11699.625
$ws.Cells.Item(69, 13).Value = -10901
$ws.Cells.Item(69, 14).Value = -13447.625
# Row 72 (Leve Item ID 12616)
$ws.Cells.Item(72, 8).Value = 3902.6667
$ws.Cells.Item(72, 9).Value = 3925
$ws.Cells.Item(72, 10).Value = 3899.875
$ws.Cells.Item(72, 11).Value = 35325
$ws.Cells.Item(72, 12).Value = 35098.875
$ws.Cells.Item(72, 13).Value = -30957
$ws.Cells.Item(72, 14).Value = -43834.875
# Row 76 (Leve Item ID 12602)
$ws.Cells.Item(76, 8).Value = 8392.857
$ws.Cells.Item(76, 9).Value = 7000
$ws.Cells.Item(76, 10).Value = 8625
$ws.Cells.Item(76, 11).Value = 7000
$ws.Cells.Item(76, 12).Value = 8625
$ws.Cells.Item(76, 13).Value = -6685
$ws.Cells.Item(76, 14).Value = -9255
# Row 79 (Leve Item ID 12602)
$ws.Cells.Item(79, 8).Value = 8392.857
$ws.Cells.Item(79, 9).Value = 7000
$ws.Cells.Item(79, 10).Value = 8625
$ws.Cells.Item(79, 11).Value = 7000
$ws.Cells.Item(79, 12).Value = 8625
$ws.Cells.Item(79, 13).Value = -5908
$ws.Cells.Item(79, 14).Value = -10809
# Row 81 (Leve Item ID 10637)
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()
# Row 84 (Leve Item ID 10637)
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()
# Row 86 (Leve Item ID 12603)
$ws.Cells.Item(86, 8).Value = 2084.8125
$ws.Cells.Item(86, 9).Value = 1101
$ws.Cells.Item(86, 10).Value = 2850
$ws.Cells.Item(86, 11).Value = 1101
$ws.Cells.Item(86, 12).Value = 2850
$ws.Cells.Item(86, 13).Value = 22
$ws.Cells.Item(86, 14).Value = -5096
# Row 89 (Leve Item ID 12603)
$ws.Cells.Item(89, 8).Value = 2084.8125
$ws.Cells.Item(89, 9).Value = 1101
$ws.Cells.Item(89, 10).Value = 2850
$ws.Cells.Item(89, 11).Value = 5505
$ws.Cells.Item(89, 12).Value = 14250
$ws.Cells.Item(89, 13).Value = 111
$ws.Cells.Item(89, 14).Value = -25482

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 63 (Leve Item ID 12528)
$ws.Cells.Item(63, 8).Value = 1550
$ws.Cells.Item(63, 9).Value = 1671.4286
$ws.Cells.Item(63, 10).Value = 1380
$ws.Cells.Item(63, 11).Value = 1671.4286
$ws.Cells.Item(63, 12).Value = 1380
$ws.Cells.Item(63, 13).Value = -985.4286
$ws.Cells.Item(63, 14).Value = -2752
# Row 66 (Leve Item ID 12528)
$ws.Cells.Item(66, 8).Value = 1550
$ws.Cells.Item(66, 9).Value = 1671.4286
$ws.Cells.Item(66, 10).Value = 1380
$ws.Cells.Item(66, 11).Value = 8357.143
$ws.Cells.Item(66, 12).Value = 6900
$ws.Cells.Item(66, 13).Value = -4925.143
$ws.Cells.Item(66, 14).Value = -13764
# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 4733.9136
$ws.Cells.Item(132, 9).Value = 3012.02
$ws.Cells.Item(132, 10).Value = 15495.75
$ws.Cells.Item(132, 11).Value = 9036.059999999999
$ws.Cells.Item(132, 12).Value = 46487.25
$ws.Cells.Item(132, 13).Value = -6506.059999999999
$ws.Cells.Item(132, 14).Value = -51547.25

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 105 (Leve Item ID 19947)
$ws.Cells.Item(105, 8).Value = 1857.2354
$ws.Cells.Item(105, 9).Value = 1378.1818
$ws.Cells.Item(105, 10).Value = 2735.5
$ws.Cells.Item(105, 11).Value = 1378.1818
$ws.Cells.Item(105, 12).Value = 2735.5
$ws.Cells.Item(105, 13).Value = 368.8181999999999
$ws.Cells.Item(105, 14).Value = -6229.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Cells.Item(22, 8).Value = 77580.16
$ws.Cells.Item(22, 9).Value = 100294.7
$ws.Cells.Item(22, 10).Value = 1865
$ws.Cells.Item(22, 11).Value = 100294.7
$ws.Cells.Item(22, 12).Value = 1865
$ws.Cells.Item(22, 13).Value = -99944.7
$ws.Cells.Item(22, 14).Value = -2565
# Row 62 (Leve Item ID 12580)
$ws.Cells.Item(62, 8).Value = 4585.185
$ws.Cells.Item(62, 9).Value = 5013.636
$ws.Cells.Item(62, 10).Value = 2700
$ws.Cells.Item(62, 11).Value = 5013.636
$ws.Cells.Item(62, 12).Value = 2700
$ws.Cells.Item(62, 13).Value = -4389.636
$ws.Cells.Item(62, 14).Value = -3948
# Row 65 (Leve Item ID 12580)
$ws.Cells.Item(65, 8).Value = 4585.185
$ws.Cells.Item(65, 9).Value = 5013.636
$ws.Cells.Item(65, 10).Value = 2700
$ws.Cells.Item(65, 11).Value = 25068.18
$ws.Cells.Item(65, 12).Value = 13500
$ws.Cells.Item(65, 13).Value = -21948.18
$ws.Cells.Item(65, 14).Value = -19740

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Cells.Item(12, 8).Value = 38.475
$ws.Cells.Item(12, 9).Value = 2.5
$ws.Cells.Item(12, 10).Value = 44.82353
$ws.Cells.Item(12, 11).Value = 7.5
$ws.Cells.Item(12, 12).Value = 134.47059
$ws.Cells.Item(12, 13).Value = 165.5
$ws.Cells.Item(12, 14).Value = -480.47059
# Row 104 (Leve Item ID 19807)
$ws.Cells.Item(104, 8).Value = 3269.5
$ws.Cells.Item(104, 10).Value = 3269.5
$ws.Cells.Item(104, 12).Value = 9808.5
$ws.Cells.Item(104, 14).Value = -15050.5
# Row 106 (Leve Item ID 19819)
$ws.Cells.Item(106, 8).Value = 38450
$ws.Cells.Item(106, 10).Value = 34600
$ws.Cells.Item(106, 12).Value = 103800
$ws.Cells.Item(106, 14).Value = -105692

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Cells.Item(70, 8).Value = 3992.3635
$ws.Cells.Item(70, 9).Value = 3984.5
$ws.Cells.Item(70, 10).Value = 4001.8
$ws.Cells.Item(70, 11).Value = 3984.5
$ws.Cells.Item(70, 12).Value = 4001.8
$ws.Cells.Item(70, 13).Value = -3714.5
$ws.Cells.Item(70, 14).Value = -4541.8
# Row 73 (Leve Item ID 14146)
$ws.Cells.Item(73, 8).Value = 3992.3635
$ws.Cells.Item(73, 9).Value = 3984.5
$ws.Cells.Item(73, 10).Value = 4001.8
$ws.Cells.Item(73, 11).Value = 3984.5
$ws.Cells.Item(73, 12).Value = 4001.8
$ws.Cells.Item(73, 13).Value = -3048.5
$ws.Cells.Item(73, 14).Value = -5873.8
# Row 80 (Leve Item ID 12521)
$ws.Cells.Item(80, 8).Value = 2643.8572
$ws.Cells.Item(80, 9).Value = 2295
$ws.Cells.Item(80, 10).Value = 2905.5
$ws.Cells.Item(80, 11).Value = 2295
$ws.Cells.Item(80, 12).Value = 2905.5
$ws.Cells.Item(80, 13).Value = -1297
$ws.Cells.Item(80, 14).Value = -4901.5
# Row 83 (Leve Item ID 12521)
$ws.Cells.Item(83, 8).Value = 2643.8572
$ws.Cells.Item(83, 9).Value = 2295
$ws.Cells.Item(83, 10).Value = 2905.5
$ws.Cells.Item(83, 11).Value = 11475
$ws.Cells.Item(83, 12).Value = 14527.5
$ws.Cells.Item(83, 13).Value = -6483
$ws.Cells.Item(83, 14).Value = -24511.5
# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 21767.102
$ws.Cells.Item(132, 9).Value = 36187.25
$ws.Cells.Item(132, 10).Value = 4676.5557
$ws.Cells.Item(132, 11).Value = 108561.75
$ws.Cells.Item(132, 12).Value = 14029.6671
$ws.Cells.Item(132, 13).Value = -106031.75
$ws.Cells.Item(132, 14).Value = -19089.6671

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Cells.Item(40, 8).Value = 2714
$ws.Cells.Item(40, 9).Value = 2550.5715
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 11).Value = 2550.5715
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 13).Value = -2414.5715
$ws.Cells.Item(40, 14).Value = -3272
# Row 68 (Leve Item ID 12563)
$ws.Cells.Item(68, 8).Value = 2401.4614
$ws.Cells.Item(68, 9).Value = 1919.8
$ws.Cells.Item(68, 10).Value = 2702.5
$ws.Cells.Item(68, 11).Value = 1919.8
$ws.Cells.Item(68, 12).Value = 2702.5
$ws.Cells.Item(68, 13).Value = -1170.8
$ws.Cells.Item(68, 14).Value = -4200.5
# Row 71 (Leve Item ID 12563)
$ws.Cells.Item(71, 8).Value = 2401.4614
$ws.Cells.Item(71, 9).Value = 1919.8
$ws.Cells.Item(71, 10).Value = 2702.5
$ws.Cells.Item(71, 11).Value = 9599
$ws.Cells.Item(71, 12).Value = 13512.5
$ws.Cells.Item(71, 13).Value = -5855
$ws.Cells.Item(71, 14).Value = -21000.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 62 (Leve Item ID 12589)
$ws.Cells.Item(62, 8).Value = 3034
$ws.Cells.Item(62, 9).Value = 2750
$ws.Cells.Item(62, 10).Value = 3223.3333
$ws.Cells.Item(62, 11).Value = 2750
$ws.Cells.Item(62, 12).Value = 3223.3333
$ws.Cells.Item(62, 13).Value = -2126
$ws.Cells.Item(62, 14).Value = -4471.3333
# Row 65 (Leve Item ID 12589)
$ws.Cells.Item(65, 8).Value = 3034
$ws.Cells.Item(65, 9).Value = 2750
$ws.Cells.Item(65, 10).Value = 3223.3333
$ws.Cells.Item(65, 11).Value = 13750
$ws.Cells.Item(65, 12).Value = 16116.6665
$ws.Cells.Item(65, 13).Value = -10630
$ws.Cells.Item(65, 14).Value = -22356.6665

